$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Don Kai's Event"
$ws.Range("A3").Select()
